$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "In Transit Qty" (column E) values
$ws.Range("E2").Value = 1009
$ws.Range("E3").Value = 509
$ws.Range("E4").Value = 1059
$ws.Range("E5").Value = 1609

# Clear the Product SKU Name in B3 (no longer populated)
$ws.Range("B3").ClearContents()

# Clear the Unrestricted Qty in D4 (no longer populated)
$ws.Range("D4").ClearContents()

# Update Expiry Date for row 4
$ws.Range("F4").Value = 43259

# Update the active selection to D4
$ws.Range("D4").Select()

$wb.Save()
